$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need to be forced
# to Text format first, otherwise Excel auto-converts them to a Number type
# (the source data stores these as text, e.g. "405.67", not numeric 405.67).
$textCells = @("D5","D6","D7","D9","D11","D14","D15","D17","D23","D24","D25","D26","D27","D28","D29","D32","D33","D34","D36","D37","D42","D43","D44","D46","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.381.11"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.395.24"
$ws.Range("E3").Value = "  +3.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value = "405.67"
$ws.Range("E5").Value = "  -0.87%  "

# Row 6 - Solana
$ws.Range("D6").Value = "129.91"
$ws.Range("E6").Value = "  +15.91%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  +7.91%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.12%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.676"
$ws.Range("E9").Value = "  +9.33%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +10.21%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "42.39"
$ws.Range("E11").Value = "  +9.66%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.42%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.948.23"
$ws.Range("E13").Value = "  +3.08%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "8.53"
$ws.Range("E14").Value = "  +4.95%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "19.73"

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.389.94"
$ws.Range("E16").Value = "  +2.47%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "11.50"
$ws.Range("E17").Value = "  +9.73%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "61.325.18"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  +4.75%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +16.84%  "

# Row 21 - ImmutableX
$ws.Range("E21").Value = "  +0.54%  "

# Row 22 - Litecoin
$ws.Range("E22").Value = "  +13.61%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "13.13"
$ws.Range("E23").Value = "  +5.67%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "308.38"
$ws.Range("E24").Value = "  +4.61%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  +2.98%  "

# Row 26 - was LEO, now Filecoin
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").Value = "8.55"
$ws.Range("E26").Value = "  +15.32%  "

# Row 27 - was Filecoin, now LEO
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "4.74"
$ws.Range("E27").Value = "  +10.38%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "29.72"
$ws.Range("E28").Value = "  +2.31%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "7.49"
$ws.Range("E29").Value = "  +2.66%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  +0.56%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +5.72%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "11.81"
$ws.Range("E32").Value = "  +6.45%  "

# Row 33 - Toncoin
$ws.Range("D33").Value = "2.62"
$ws.Range("E33").Value = "  +7.05%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "42.63"

# Row 36 - VeChain
$ws.Range("D36").Value = "0.0488"
$ws.Range("E36").Value = "  +2.45%  "

# Row 37 - OKB
$ws.Range("D37").Value = "52.29"
$ws.Range("E37").Value = "  -0.03%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.24%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +5.33%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -0.17%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +8.75%  "

# Row 42 - Stellar
$ws.Range("D42").Value = "0.125"
$ws.Range("E42").Value = "  +4.92%  "

# Row 43 - Monero
$ws.Range("D43").Value = "137.13"
$ws.Range("E43").Value = "  +1.78%  "

# Row 44 - NEARProtocol
$ws.Range("D44").Value = "4.03"
$ws.Range("E44").Value = "  +7.83%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  +0.83%  "

# Row 46 - Celestia
$ws.Range("D46").Value = "17.06"
$ws.Range("E46").Value = "  +5.45%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  +1.38%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "21.75"
$ws.Range("E48").Value = "  +5.04%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.150.64"
$ws.Range("E49").Value = "  +2.22%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "3.731.87"
$ws.Range("E50").Value = "  +2.46%  "

# Row 51 - ApeXProtocol
$ws.Range("E51").Value = "  +0.22%  "
